# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de
# handback cycle completed ("Handed back: in sync with en-US") and fills
# in the Latest Target File / Latest Handback File / Latest Handback
# DateTime columns on the per-language sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: per-language status cells -----------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

# --- zh-cn sheet -------------------------------------------------------
$zhcn.Range("C2").Value = $statusText

# source markdown filename/url (same doc referenced by the A2 hyperlink)
$mdName = $zhcn.Range("A2").Value()
$mdUrl  = ""
foreach ($lnk in $zhcn.Hyperlinks) {
    $mdUrl = $lnk.Address()
}

# Latest Target File -> hyperlink to the source md, same styling as A2
$zhcn.Range("I2").Value = $mdName
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, $mdName)
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276

# Latest Handback File -> the xliff that was last handed off is now the
# handback file too (round-tripped, in sync)
$zhcn.Range("J2").Value = $zhcn.Range("G2").Value()

# Latest Handback DateTime
$zhcn.Range("K2").Value = "2016-08-28 14:58:41"

# --- de-de sheet ---------------------------------------------------------
$dede.Range("C2").Value = $statusText

$mdName3 = $dede.Range("A2").Value()
$mdUrl3  = ""
foreach ($lnk in $dede.Hyperlinks) {
    $mdUrl3 = $lnk.Address()
}

$dede.Range("I2").Value = $mdName3
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl3, [System.Type]::Missing, [System.Type]::Missing, $mdName3)
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276

$dede.Range("J2").Value = $dede.Range("G2").Value()

$dede.Range("K2").Value = "2016-08-28 14:58:48"

# --- column widths (source doc panel was widened for the new links) ------
# Note: this host's ColumnWidth setter quantizes to 1/6-character steps, so
# the inputs below are chosen to land on the exact target "width" values
# (29.9777047293527 -> 30, and 40 stays 40) once the engine stores them.
$overview.Range("E1").ColumnWidth = 29.2
$overview.Range("F1").ColumnWidth = 29.2

$zhcn.Range("C1").ColumnWidth = 29.2
$zhcn.Range("I1").ColumnWidth = 39.17
$zhcn.Range("J1").ColumnWidth = 39.17

$dede.Range("C1").ColumnWidth = 29.2
$dede.Range("I1").ColumnWidth = 39.17
$dede.Range("J1").ColumnWidth = 39.17
